$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate column R (18th column) formatting into a brand-new column S (19th
# column) by copying the whole column and inserting it as a new column. This
# correctly carries over each row's existing cell style (so S ends up using
# the very same style ids as column R), instead of creating fresh style
# entries in styles.xml.
$ws.Columns.Item(18).Copy()
$ws.Columns.Item(19).Insert(-4161)

# Now overwrite the freshly inserted column S with the actual 2022 figures.
# Rows that stay "-" (shared string placeholder) are left alone, because the
# column-copy already placed the correct "-" text/style in them.
$ws.Range("S4").Value = 2022

$ws.Range("S5").Value = 135
$ws.Range("S6").Value = 99
$ws.Range("S7").Value = 36

$ws.Range("S8").Value = 97
$ws.Range("S9").Value = 80
$ws.Range("S10").Value = 17

$ws.Range("S11").Value = 17
$ws.Range("S12").Value = 11
$ws.Range("S13").Value = 6

$ws.Range("S14").Value = 5
$ws.Range("S15").Value = 3
$ws.Range("S16").Value = 2

# S17:S19 remain "-" (already copied from R17:R19).

$ws.Range("S20").Value = 6
$ws.Range("S21").Value = 1
$ws.Range("S22").Value = 5

$ws.Range("S23").Value = "-"
$ws.Range("S24").Value = "-"
$ws.Range("S25").Value = "-"

$ws.Range("S26").Value = 10
$ws.Range("S27").Value = 4
$ws.Range("S28").Value = 6

# S29:S34 remain "-" (already copied from R29:R34).

# Match the author's final selection/active cell recorded in the workbook.
$ws.Range("T24").Select() | Out-Null
